$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.860.72'
$ws.Range("E2").Value = '  -0.92%  '
$ws.Range("D3").Value = '1.562.94'
$ws.Range("E3").Value = '  +0.14%  '
$ws.Range("E4").Value = '  -0.07%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '205.89'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.32%  '
$ws.Range("E7").Value = '  -0.02%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '21.77'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -1.50%  '
$ws.Range("E9").Value = '  -0.17%  '
$ws.Range("E10").Value = '  -0.71%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0865'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -0.05%  '
$ws.Range("D12").Value = '1.785.11'
$ws.Range("E12").Value = '  +0.15%  '
$ws.Range("D13").Value = '1.567.31'
$ws.Range("E13").Value = '  +0.51%  '
$ws.Range("D16").Value = '26.875.04'
$ws.Range("E16").Value = '  -0.79%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '61.32'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -2.54%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '215.54'
$ws.Range("D18").Style = "Normal"
$ws.Range("E19").Value = '  +2.13%  '
$ws.Range("E20").Value = '  -0.26%  '
$ws.Range("E21").Value = '  -0.02%  '
$ws.Range("E22").Value = '  +0.52%  '
$ws.Range("E23").Value = '  -1.19%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.02'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +1.26%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '154.07'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +1.57%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '6.69'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +1.37%  '
$ws.Range("E27").Value = '  +0.33%  '
$ws.Range("E28").Value = '  -0.06%  '
$ws.Range("E29").Value = '  -0.83%  '
$ws.Range("E30").Value = '  +0.86%  '
$ws.Range("E31").Value = '  -3.48%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.17'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +0.32%  '
$ws.Range("D33").Value = '1.395.61'
$ws.Range("E33").Value = '  +0.87%  '
$ws.Range("E34").Value = '  -0.22%  '
$ws.Range("E35").Value = '  -1.23%  '
$ws.Range("E36").Value = '  -0.52%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.919'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -2.78%  '
$ws.Range("E38").Value = '  -0.50%  '
$ws.Range("E39").Value = '  +2.86%  '
$ws.Range("E40").Value = '  +0.31%  '
$ws.Range("E41").Value = '  -0.06%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.990'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +0.15%  '
$ws.Range("E44").Value = '  -0.18%  '
$ws.Range("D47").Value = '1.699.01'
$ws.Range("E47").Value = '  +0.27%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '86.70'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +1.72%  '
$ws.Range("E49").Value = '  +2.76%  '
$ws.Range("D50").Value = '0.0₇0982'
$ws.Range("E50").Value = '  +0.05%  '
$ws.Range("E51").Value = '  +1.03%  '
